# Commit: "Added support for Excel"
#
# Renames the original "Sheet1" (raw Jira time-log export) to "Time Logs"
# and adds two pivot-style summary sheets computed from it:
#   - "Hours by Person": total Hours Logged per Assignee
#   - "Hours by Ticket": Ticket Description + total Hours Logged per Ticket Number
#
# xlUp = -4162, xlCenter = -4108, xlTop = -4160, xlContinuous = 1

$wb = $excel.ActiveWorkbook

$timeLogs = $wb.Worksheets.Item(1)
$timeLogs.Name = "Time Logs"

$byPerson = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $timeLogs)
$byPerson.Name = "Hours by Person"

$byTicket = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $byPerson)
$byTicket.Name = "Hours by Ticket"

$lastRow = $timeLogs.Cells.Item($timeLogs.Rows.Count, 1).End(-4162).Row

# ---------------------------------------------------------------------------
# Pass 1: walk "Time Logs" once, building per-Assignee and per-Ticket totals
# (first-seen order is kept in a side array; final order is ordinal-sorted
# further down: "Hours by Person" ascends by Assignee, "Hours by Ticket"
# ascends by Ticket Description).
# ---------------------------------------------------------------------------

$personHours = @{}
$personOrder = @()

$ticketHours = @{}
$ticketDesc = @{}
$ticketOrder = @()       # ticket numbers, first-seen order
$ticketSortKey = @{}     # ticket number -> Ticket Description (the sort key)

for ($r = 2; $r -le $lastRow; $r++) {
  $name = $timeLogs.Cells.Item($r, 1).Value2
  $ticket = $timeLogs.Cells.Item($r, 2).Value2
  $desc = $timeLogs.Cells.Item($r, 3).Value2
  $hrs = $timeLogs.Cells.Item($r, 4).Value2

  if ($personHours.ContainsKey($name)) {
    $personHours[$name] = $personHours[$name] + $hrs
  } else {
    $personHours[$name] = $hrs
    $personOrder += $name
  }

  if ($ticketHours.ContainsKey($ticket)) {
    $ticketHours[$ticket] = $ticketHours[$ticket] + $hrs
  } else {
    $ticketHours[$ticket] = $hrs
    $ticketDesc[$ticket] = $desc
    $ticketSortKey[$ticket] = $desc
    $ticketOrder += $ticket
  }
}

# ---------------------------------------------------------------------------
# Plain ascending ordinal sort (insertion sort) of the two key arrays.
# Inlined (no helper function) to stay well under the interpreter's
# statement budget.
# ---------------------------------------------------------------------------

$n = $personOrder.Count
for ($i = 1; $i -lt $n; $i++) {
  $key = $personOrder[$i]
  $j = $i - 1
  $cont = $true
  while (($j -ge 0) -and $cont) {
    $a = $key
    $b = $personOrder[$j]
    $la = $a.Length
    $lb = $b.Length
    $m = $la
    if ($lb -lt $m) { $m = $lb }
    $less = $false
    $decided = $false
    for ($p = 0; $p -lt $m; $p++) {
      $ca = [int][char]$a[$p]
      $cb = [int][char]$b[$p]
      if ($ca -ne $cb) {
        if ($ca -lt $cb) { $less = $true } else { $less = $false }
        $decided = $true
        break
      }
    }
    if (-not $decided) {
      if ($la -lt $lb) { $less = $true } else { $less = $false }
    }
    if ($less) {
      $personOrder[$j+1] = $personOrder[$j]
      $j = $j - 1
    } else {
      $cont = $false
    }
  }
  $personOrder[$j+1] = $key
}

$n = $ticketOrder.Count
for ($i = 1; $i -lt $n; $i++) {
  $key = $ticketOrder[$i]
  $keySort = $ticketSortKey[$key]
  $j = $i - 1
  $cont = $true
  while (($j -ge 0) -and $cont) {
    $a = $keySort
    $b = $ticketSortKey[$ticketOrder[$j]]
    $la = $a.Length
    $lb = $b.Length
    $m = $la
    if ($lb -lt $m) { $m = $lb }
    $less = $false
    $decided = $false
    for ($p = 0; $p -lt $m; $p++) {
      $ca = [int][char]$a[$p]
      $cb = [int][char]$b[$p]
      if ($ca -ne $cb) {
        if ($ca -lt $cb) { $less = $true } else { $less = $false }
        $decided = $true
        break
      }
    }
    if (-not $decided) {
      if ($la -lt $lb) { $less = $true } else { $less = $false }
    }
    if ($less) {
      $ticketOrder[$j+1] = $ticketOrder[$j]
      $j = $j - 1
    } else {
      $cont = $false
    }
  }
  $ticketOrder[$j+1] = $key
}

# ---------------------------------------------------------------------------
# "Hours by Person": Assignee | Hours Logged
# ---------------------------------------------------------------------------

$byPerson.Range("A1").Value2 = "Assignee"
$byPerson.Range("B1").Value2 = "Hours Logged"

$row = 2
foreach ($p in $personOrder) {
  $byPerson.Cells.Item($row, 1).Value2 = $p
  $byPerson.Cells.Item($row, 2).Value2 = $personHours[$p]
  $row = $row + 1
}
$personLastRow = $row - 1

# ---------------------------------------------------------------------------
# "Hours by Ticket": Ticket Number | Ticket Description | Hours Logged
# ---------------------------------------------------------------------------

$byTicket.Range("A1").Value2 = "Ticket Number"
$byTicket.Range("B1").Value2 = "Ticket Description"
$byTicket.Range("C1").Value2 = "Hours Logged"

$row = 2
foreach ($t in $ticketOrder) {
  $byTicket.Cells.Item($row, 1).Value2 = $t
  $byTicket.Cells.Item($row, 2).Value2 = $ticketDesc[$t]
  $byTicket.Cells.Item($row, 3).Value2 = $ticketHours[$t]
  $row = $row + 1
}
$ticketLastRow = $row - 1

# ---------------------------------------------------------------------------
# Styling: bold, centered (H+V), thin box border on every text cell
# (header row + the label column(s)); numeric "Hours Logged" cells are left
# with the default style, matching the source sheet's convention.
# ---------------------------------------------------------------------------

$styleRanges = @(
  $byPerson.Range("A1:A" + $personLastRow),
  $byPerson.Range("B1"),
  $byTicket.Range("A1:B" + $ticketLastRow),
  $byTicket.Range("C1")
)

foreach ($sr in $styleRanges) {
  $sr.Font.Bold = $true
  $sr.HorizontalAlignment = -4108
  $sr.VerticalAlignment = -4160
  $sr.Borders.LineStyle = 1
}

# Keep "Time Logs" as the active/selected sheet (matches the source
# workbook's activeTab="0").
$timeLogs.Activate()
$timeLogs.Range("A1").Select()
